$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 66.334
$ws.Range("D2").Value = 66.334
$ws.Range("E2").Value = 2.61336138
$ws.Range("F2").Value = 0.007458129999999999
$ws.Range("G2").Value = 0.48770223
$ws.Range("H2").Value = 32.71349691
$ws.Range("I2").Value = 10.59484283482746
$ws.Range("J2").Value = 10.59484283482746
$ws.Range("K2").Value = 0.4200334952196099
$ws.Range("L2").Value = 0.000961669283543418
$ws.Range("M2").Value = 0.05544099550973278
$ws.Range("N2").Value = 8.118351175721624

# Row 3
$ws.Range("C3").Value = 111.812
$ws.Range("D3").Value = 111.812
$ws.Range("E3").Value = 1.56395727
$ws.Range("F3").Value = 0.00327679
$ws.Range("G3").Value = 0.3584990199999999
$ws.Range("H3").Value = 40.57860195
$ws.Range("I3").Value = 21.30061906077784
$ws.Range("J3").Value = 21.30061906077784
$ws.Range("K3").Value = 0.2855045774280677
$ws.Range("L3").Value = 0.0008370201901380269
$ws.Range("M3").Value = 0.08407043989121922
$ws.Range("N3").Value = 13.9042836686778

# Row 4
$ws.Range("C4").Value = 44.876
$ws.Range("D4").Value = 89.711
$ws.Range("E4").Value = 1.98572922
$ws.Range("F4").Value = 0.01286585
$ws.Range("G4").Value = 0.28449185
$ws.Range("H4").Value = 13.22304504
$ws.Range("I4").Value = 9.979690286688172
$ws.Range("J4").Value = 19.95855262344499
$ws.Range("K4").Value = 0.4779556440546291
$ws.Range("L4").Value = 0.002083452291498644
$ws.Range("M4").Value = 0.0625861210735138
$ws.Range("N4").Value = 5.44044955023412

# Row 5
$ws.Range("C5").Value = 64.01300000000001
$ws.Range("D5").Value = 124.812
$ws.Range("E5").Value = 1.39916202
$ws.Range("F5").Value = 0.0061942
$ws.Range("G5").Value = 0.19510707
$ws.Range("H5").Value = 12.74446728
$ws.Range("I5").Value = 12.43710225544972
$ws.Range("J5").Value = 22.62429036488249
$ws.Range("K5").Value = 0.2565288253146943
$ws.Range("L5").Value = 0.001290392061224677
$ws.Range("M5").Value = 0.04281632394673101
$ws.Range("N5").Value = 4.573077103208242

# Row 6
$ws.Range("C6").Value = 26.103
$ws.Range("D6").Value = 104.304
$ws.Range("E6").Value = 1.75244385
$ws.Range("F6").Value = 0.02182721
$ws.Range("G6").Value = 0.14077328
$ws.Range("H6").Value = 3.90214375
$ws.Range("I6").Value = 7.024794569148534
$ws.Range("J6").Value = 28.08837020085965
$ws.Range("K6").Value = 0.5260207542958999
$ws.Range("L6").Value = 0.004014096247548704
$ws.Range("M6").Value = 0.04027511750860489
$ws.Range("N6").Value = 2.036979954913046

# Row 7
$ws.Range("C7").Value = 33.575
$ws.Range("D7").Value = 124.614
$ws.Range("E7").Value = 1.41354863
$ws.Range("F7").Value = 0.01192575
$ws.Range("G7").Value = 0.0986142
$ws.Range("H7").Value = 3.42103651
$ws.Range("I7").Value = 7.397776835110809
$ws.Range("J7").Value = 25.31161159787503
$ws.Range("K7").Value = 0.2954796991036872
$ws.Range("L7").Value = 0.002713180677154931
$ws.Range("M7").Value = 0.02567685665590178
$ws.Range("N7").Value = 1.519794786646947

# Row 8
$ws.Range("C8").Value = 17.301
$ws.Range("D8").Value = 103.699
$ws.Range("E8").Value = 1.80872041
$ws.Range("F8").Value = 0.03037105
$ws.Range("G8").Value = 0.08652232000000001
$ws.Range("H8").Value = 1.62807922
$ws.Range("I8").Value = 5.442243748401665
$ws.Range("J8").Value = 32.63524871145398
$ws.Range("K8").Value = 0.613659465247151
$ws.Range("L8").Value = 0.006752611299109171
$ws.Range("M8").Value = 0.02998056304638354
$ws.Range("N8").Value = 1.019980528351894

# Row 9
$ws.Range("C9").Value = 22.043
$ws.Range("D9").Value = 113.324
$ws.Range("E9").Value = 1.55821096
$ws.Range("F9").Value = 0.01586193
$ws.Range("G9").Value = 0.05736815
$ws.Range("H9").Value = 1.31241925
$ws.Range("I9").Value = 5.069995646532239
$ws.Range("J9").Value = 23.5688462703314
$ws.Range("K9").Value = 0.3349562838341122
$ws.Range("L9").Value = 0.003649305496352108
$ws.Range("M9").Value = 0.01552996932465498
$ws.Range("N9").Value = 0.6119261487629675

# Row 10
$ws.Range("C10").Value = 12.461
$ws.Range("D10").Value = 99.51000000000001
$ws.Range("E10").Value = 1.91079858
$ws.Range("F10").Value = 0.03441491000000001
$ws.Range("G10").Value = 0.05320871
$ws.Range("H10").Value = 0.7386829799999999
$ws.Range("I10").Value = 4.290324683038089
$ws.Range("J10").Value = 34.28038863614363
$ws.Range("K10").Value = 0.6894536043448665
$ws.Range("L10").Value = 0.00737419948264962
$ws.Range("M10").Value = 0.02058499628146208
$ws.Range("N10").Value = 0.5517653478985983

# Row 11
$ws.Range("C11").Value = 16.209
$ws.Range("D11").Value = 100.493
$ws.Range("E11").Value = 1.75045176
$ws.Range("F11").Value = 0.01833839
$ws.Range("G11").Value = 0.03641956
$ws.Range("H11").Value = 0.61289786
$ws.Range("I11").Value = 3.988888345169645
$ws.Range("J11").Value = 20.69807627199594
$ws.Range("K11").Value = 0.3487761036105429
$ws.Range("L11").Value = 0.004521547223815891
$ws.Range("M11").Value = 0.01022278953281502
$ws.Range("N11").Value = 0.3008847876647613

# Row 12
$ws.Range("C12").Value = 9.282
$ws.Range("D12").Value = 92.58499999999999
$ws.Range("E12").Value = 2.06003375
$ws.Range("F12").Value = 0.03969014
$ws.Range("G12").Value = 0.03658754
$ws.Range("H12").Value = 0.37911152
$ws.Range("I12").Value = 3.282567783192184
$ws.Range("J12").Value = 32.76263367260277
$ws.Range("K12").Value = 0.7347914006938098
$ws.Range("L12").Value = 0.008909707720911616
$ws.Range("M12").Value = 0.01425929407846041
$ws.Range("N12").Value = 0.2792872139837131

# Row 13
$ws.Range("C13").Value = 12.936
$ws.Range("D13").Value = 89.84
$ws.Range("E13").Value = 1.96628306
$ws.Range("F13").Value = 0.02040392
$ws.Range("G13").Value = 0.02591663
$ws.Range("H13").Value = 0.35126333
$ws.Range("I13").Value = 3.387828416756282
$ws.Range("J13").Value = 19.43852505351824
$ws.Range("K13").Value = 0.4124571842915458
$ws.Range("L13").Value = 0.005501379718980723
$ws.Range("M13").Value = 0.008305930436954281
$ws.Range("N13").Value = 0.1936688846565064
